$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (price & volume refresh, plus a couple of row swaps)

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '25.889.74'
$ws.Cells.Item(2, 5).Value = '  +0.51%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.744.29'
$ws.Cells.Item(3, 5).Value = '  -0.54%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9987'
$ws.Cells.Item(4, 5).Value = '  -0.30%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '226.07'
$ws.Cells.Item(5, 5).Value = '  -4.69%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9992'
$ws.Cells.Item(6, 5).Value = '  -0.19%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5153'
$ws.Cells.Item(7, 5).Value = '  +1.67%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2802'
$ws.Cells.Item(8, 5).Value = '  +6.65%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '39.19'
$ws.Cells.Item(9, 5).Value = '  -3.46%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.06095'
$ws.Cells.Item(10, 5).Value = '  -1.85%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.745.48'
$ws.Cells.Item(11, 5).Value = '  -0.51%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.06958'
$ws.Cells.Item(12, 5).Value = '  +0.14%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '15.24'
$ws.Cells.Item(13, 5).Value = '  -1.60%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6335'
$ws.Cells.Item(14, 5).Value = '  +4.78%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.482'
$ws.Cells.Item(15, 5).Value = '  +0.66%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '76.45'
$ws.Cells.Item(16, 5).Value = '  -2.47%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.9973'
$ws.Cells.Item(17, 5).Value = '  -0.36%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.9982'
$ws.Cells.Item(18, 5).Value = '  -0.22%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '25.891.26'
$ws.Cells.Item(19, 5).Value = '  +0.32%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '11.44'
$ws.Cells.Item(20, 5).Value = '  -1.89%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.000006586'
$ws.Cells.Item(21, 5).Value = '  -2.73%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '1.966.57'
$ws.Cells.Item(22, 5).Value = '  -0.40%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '4.090'
$ws.Cells.Item(23, 5).Value = '  +0.77%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '8.409'
$ws.Cells.Item(24, 5).Value = '  +2.75%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '5.118'
$ws.Cells.Item(25, 5).Value = '  -1.18%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '138.26'
$ws.Cells.Item(26, 5).Value = '  +0.26%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '1.514'
$ws.Cells.Item(27, 5).Value = '  +3.06%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.825'
$ws.Cells.Item(28, 5).Value = '  +1.11%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '15.00'
$ws.Cells.Item(29, 5).Value = '  -0.56%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '102.93'
$ws.Cells.Item(30, 5).Value = '  +0.25%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08300'
$ws.Cells.Item(31, 5).Value = '  +0.28%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.614'
$ws.Cells.Item(32, 5).Value = '  -2.35%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.404'
$ws.Cells.Item(33, 5).Value = '  +0.20%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.04381'
$ws.Cells.Item(34, 5).Value = '  +0.06%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.95%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.9682'
$ws.Cells.Item(36, 5).Value = '  -3.40%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.6049'
$ws.Cells.Item(37, 5).Value = '  +0.63%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.668'
$ws.Cells.Item(38, 5).Value = '  -1.09%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01558'
$ws.Cells.Item(39, 5).Value = '  +0.55%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.906'
$ws.Cells.Item(40, 5).Value = '  -2.83%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'PaxDollar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.9983'
$ws.Cells.Item(41, 5).Value = '  -0.25%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '100.48'
$ws.Cells.Item(42, 5).Value = '  -2.94%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.3823'
$ws.Cells.Item(43, 5).Value = '  +0.30%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.7192'
$ws.Cells.Item(44, 5).Value = '  -3.89%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '4.926'
$ws.Cells.Item(45, 5).Value = '  +0.84%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.05440'
$ws.Cells.Item(46, 5).Value = '  -0.82%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '6.272'
$ws.Cells.Item(47, 5).Value = '  +5.34%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.1103'
$ws.Cells.Item(48, 5).Value = '  +2.16%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '52.33'
$ws.Cells.Item(49, 5).Value = '  +0.57%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '29.65'
$ws.Cells.Item(50, 5).Value = '  -1.83%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '7.486'
$ws.Cells.Item(51, 5).Value = '  +0.19%  '
